$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.151.93"
$ws.Range("E2").Value = "  -4.74%  "
$ws.Range("D3").Value = "3.295.92"
$ws.Range("E3").Value = "  -6.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'177.61"
$ws.Range("E5").Value = "  -11.69%  "
$ws.Range("D6").Value = "'525.60"
$ws.Range("E6").Value = "  -5.08%  "
$ws.Range("D7").Value = "'0.604"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("D8").Value = "3.288.44"
$ws.Range("E8").Value = "  -6.09%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'0.608"
$ws.Range("E10").Value = "  -7.43%  "
$ws.Range("D11").Value = "'57.37"
$ws.Range("E11").Value = "  -8.31%  "
$ws.Range("E12").Value = "  -7.55%  "
$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("D14").Value = "'9.09"
$ws.Range("E14").Value = "  -7.74%  "
$ws.Range("D15").Value = "3.814.18"
$ws.Range("E15").Value = "  -6.22%  "
$ws.Range("E16").Value = "  -5.39%  "
$ws.Range("D17").Value = "3.291.51"
$ws.Range("E17").Value = "  -6.20%  "
$ws.Range("D18").Value = "64.053.07"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "'17.45"
$ws.Range("E19").Value = "  -5.55%  "
$ws.Range("D20").Value = "'11.08"
$ws.Range("E20").Value = "  -6.36%  "
$ws.Range("D21").Value = "'0.955"
$ws.Range("E21").Value = "  -7.17%  "
$ws.Range("D22").Value = "'374.13"
$ws.Range("E22").Value = "  -4.67%  "
$ws.Range("D23").Value = "'3.78"
$ws.Range("E23").Value = "  -5.62%  "
$ws.Range("D24").Value = "'80.55"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E25").Value = "  -11.10%  "
$ws.Range("D26").Value = "'3.88"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "'6.09"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "'2.67"
$ws.Range("E28").Value = "  -5.35%  "
$ws.Range("D29").Value = "'11.38"
$ws.Range("E29").Value = "  -7.32%  "
$ws.Range("D30").Value = "'8.35"
$ws.Range("E30").Value = "  -5.86%  "
$ws.Range("D31").Value = "'28.84"
$ws.Range("E31").Value = "  -7.09%  "
$ws.Range("D32").Value = "'638.73"
$ws.Range("E32").Value = "  -7.44%  "
$ws.Range("D33").Value = "'6.63"
$ws.Range("E33").Value = "  -6.65%  "
$ws.Range("D34").Value = "'11.24"
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("D35").Value = "'59.15"
$ws.Range("E35").Value = "  -7.25%  "
$ws.Range("D36").Value = "'0.105"
$ws.Range("E36").Value = "  -5.95%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "'0.389"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").Value = "'36.62"
$ws.Range("E39").Value = "  -5.62%  "
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("D42").Value = "2.936.41"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").Value = "'2.47"
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("D45").Value = "'2.69"
$ws.Range("E45").Value = "  -10.68%  "
$ws.Range("D46").Value = "'0.0398"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("D48").Value = "'3.00"
$ws.Range("E48").Value = "  +5.08%  "
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "'135.63"
$ws.Range("E51").Value = "  -1.91%  "
